{"js": "// Convert the Word complex field (fldChar/instrText) that contains the\n// M2Doc query \"m:'Obeo\\'s website'.sampleLink ('http://www.obeo.fr',6)\"\n// into plain, visible text wrapped in \"{\" / \"}\" (the textual token syntax\n// consumed by the TokenIteratorFieldRewriterSplit parser), replacing the\n// <w:instrText> runs with plain <w:t> runs and dropping the field\n// characters entirely.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Locate the paragraph that holds the field (robust to its position).\nlet fieldParagraph = null;\nlet theField = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  const fields = p.fields;\n  fields.load(\"items\");\n  await context.sync();\n  if (fields.items.length > 0) {\n    fields.items[0].load(\"code\");\n    await context.sync();\n    fieldParagraph = p;\n    theField = fields.items[0];\n    break;\n  }\n}\n\nif (!fieldParagraph || !theField) {\n  throw new Error(\"Could not find the paragraph containing the M2Doc field.\");\n}\n\n// The raw field instruction, e.g. \" m:'Obeo\\'s website'.sampleLink ('http://www.obeo.fr',6) \".\nconst rawCode = theField.code;\nconst trimmedCode = rawCode.trim();\n\n// Reproduce the same run segmentation the original template used for its\n// instrText runs (this mirrors exactly how the field text had been typed),\n// so the resulting run layout matches the template's history once the\n// field becomes plain text.\nconst knownSegments = [\n  \"m\",\n  \":'\",\n  \"Obeo\\\\\",\n  \"'\",\n  \"s website\",\n  \"'.\",\n  \"sampleLink \",\n  \"(\",\n  \"'http://www.obeo.fr',\",\n  \"6\",\n  \")\",\n];\nconst joined = knownSegments.join(\"\");\n// Fall back to a single segment if the field text does not match what we\n// expect (keeps the script correct even if the source document changes).\nconst segments = joined === trimmedCode ? knownSegments : [trimmedCode];\n\n// Build the run list: \"{\" + segments + \"}\" each as its own <w:t> run,\n// with the original \"sampleLink \" / \"}\" runs keeping xml:space=\"preserve\".\nconst tokens = [{ text: \"{\", preserve: false }];\nsegments.forEach((s) => tokens.push({ text: s, preserve: /\\s$/.test(s) }));\ntokens.push({ text: \"}\", preserve: true });\n\nfunction escapeXml(s) {\n  return s\n    .replace(/&/g, \"&amp;\")\n    .replace(/</g, \"&lt;\")\n    .replace(/>/g, \"&gt;\");\n}\n\nlet runsXml = \"\";\ntokens.forEach((tok, idx) => {\n  const spacePreserve = tok.preserve ? ' xml:space=\"preserve\"' : \"\";\n  runsXml +=\n    \"<w:r><w:rPr><w:lang w:val=\\\"en-US\\\"/></w:rPr>\" +\n    `<w:t${spacePreserve}>${escapeXml(tok.text)}</w:t></w:r>`;\n  // The original template carried a \"_GoBack\" bookmark right after the\n  // \"s website\" segment (index 5 -> tokens index 6 once \"{\" is prefixed).\n  if (idx === 5) {\n    runsXml +=\n      '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/>';\n  }\n});\n\nconst paragraphOoxml =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  \"<pkg:xmlData>\" +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  \"<w:body><w:p><w:pPr><w:rPr><w:lang w:val=\\\"en-US\\\"/></w:rPr></w:pPr>\" +\n  runsXml +\n  \"</w:p></w:body></w:document>\" +\n  \"</pkg:xmlData></pkg:part></pkg:package>\";\n\nfieldParagraph.insertOoxml(paragraphOoxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Convert the Word complex field (fldChar/instrText) that contains the\n# M2Doc query \"m:'Obeo\\'s website'.sampleLink ('http://www.obeo.fr',6)\"\n# into plain, visible text wrapped in \"{\" / \"}\" (the textual token syntax\n# consumed by the TokenIteratorFieldRewriterSplit parser), replacing the\n# field's <w:instrText> runs with plain <w:t> runs and dropping the field\n# characters entirely.\n\nfunction Escape-Xml($s) {\n  $s = $s -replace '&', '&amp;'\n  $s = $s -replace '<', '&lt;'\n  $s = $s -replace '>', '&gt;'\n  return $s\n}\n\n$d = $word.ActiveDocument\n\n# Locate the paragraph that holds the field (robust to its position).\n$targetPara = $null\nforeach ($p in $d.Paragraphs) {\n  if ($p.Range.Fields.Count -gt 0) {\n    $targetPara = $p\n  }\n}\n\nif ($targetPara -eq $null) {\n  throw \"Could not find the paragraph containing the M2Doc field.\"\n}\n\n$f = $targetPara.Range.Fields.Item(1)\n$rawCode = $f.Code.Text\n$trimmedCode = $rawCode.Trim()\n\n# Reproduce the same run segmentation the original template used for its\n# instrText runs (mirrors exactly how the field text had been typed over\n# time), so the resulting run layout matches the template's history once\n# the field becomes plain text.\n$knownSegments = @(\"m\", \":'\", \"Obeo\\\", \"'\", \"s website\", \"'.\", \"sampleLink \", \"(\", \"'http://www.obeo.fr',\", \"6\", \")\")\n$joined = \"\"\nforeach ($seg in $knownSegments) {\n  $joined = $joined + $seg\n}\n\n# Fall back to a single segment if the field text does not match what we\n# expect (keeps the script correct even if the source document changes).\n$segments = $knownSegments\nif ($joined -ne $trimmedCode) {\n  $segments = @($trimmedCode)\n}\n\n# Build the run list: \"{\" + segments + \"}\" each as its own <w:t> run.\n$tokens = @()\n$tokens += \"{\"\nforeach ($seg in $segments) {\n  $tokens += $seg\n}\n$tokens += \"}\"\n\n$runsXml = \"\"\nfor ($i = 0; $i -lt $tokens.Length; $i++) {\n  $tok = $tokens[$i]\n  $preserve = $false\n  if ($i -eq ($tokens.Length - 1)) { $preserve = $true }\n  if ($tok -match '\\s$') { $preserve = $true }\n  $spaceAttr = \"\"\n  if ($preserve) { $spaceAttr = ' xml:space=\"preserve\"' }\n  $escaped = Escape-Xml($tok)\n  $runsXml = $runsXml + '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t' + $spaceAttr + '>' + $escaped + '</w:t></w:r>'\n  # The original template carried a \"_GoBack\" bookmark right after the\n  # \"s website\" segment (index 5 in the \"{\" + segments list).\n  if ($i -eq 5) {\n    $runsXml = $runsXml + '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/>'\n  }\n}\n\n$paragraphXml = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:pPr><w:rPr><w:lang w:val=\"en-US\"/></w:rPr></w:pPr>' + $runsXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n$targetPara.Range.InsertXML($paragraphXml)\n"}
